$wb = $excel.ActiveWorkbook

# --- Sheet1: ChosenTreatment ---
$ws1 = $wb.Worksheets.Item("ChosenTreatment")
$ws1.Range("A2").Value = "{'SexualOrientation': 'Straight or heterosexual'}"
$ws1.Range("B2").Value = "{'DevType': 'Back-end developer'}"

# --- Sheet2: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("A2").Value = 79

# --- Sheet3: Subgroups ---
$ws3 = $wb.Worksheets.Item("Subgroups")
$ws3.Cells.Item(2, 1).Value = "{'Hobby': np.int64(1), 'HDI': np.int64(1)}"
$ws3.Cells.Item(2, 2).Value = 22660
$ws3.Cells.Item(2, 3).Value = 7025.684796894849
$ws3.Cells.Item(2, 4).Value = 312.3290956815272
$ws3.Cells.Item(3, 1).Value = "{'GINI': np.int64(2), 'Hobby': np.int64(1)}"
$ws3.Cells.Item(3, 2).Value = 15435
$ws3.Cells.Item(3, 3).Value = 7774.367462050464
$ws3.Cells.Item(3, 4).Value = 1061.011760837143
$ws3.Cells.Item(4, 1).Value = "{'UndergradMajor': np.int64(2), 'Student': np.int64(1)}"
$ws3.Cells.Item(4, 2).Value = 18192
$ws3.Cells.Item(4, 3).Value = 2404.812213119229
$ws3.Cells.Item(4, 4).Value = -4308.543488094092
$ws3.Cells.Item(5, 1).Value = "{'HoursComputer': np.int64(2), 'Student': np.int64(1)}"
$ws3.Cells.Item(5, 2).Value = 16690
$ws3.Cells.Item(5, 3).Value = 4361.051980132219
$ws3.Cells.Item(5, 4).Value = -2352.303721081103
$ws3.Cells.Item(6, 1).Value = "{'Gender': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(6, 2).Value = 27820
$ws3.Cells.Item(6, 3).Value = 4860.761884746209
$ws3.Cells.Item(6, 4).Value = -1852.593816467112
$ws3.Cells.Item(7, 1).Value = "{'RaceEthnicity': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(7, 2).Value = 23264
$ws3.Cells.Item(7, 3).Value = 5618.412897795355
$ws3.Cells.Item(7, 4).Value = -1094.942803417966
$ws3.Cells.Item(8, 1).Value = "{'Dependents': np.int64(2), 'Student': np.int64(1)}"
$ws3.Cells.Item(8, 2).Value = 19464
$ws3.Cells.Item(8, 3).Value = 5774.123848787556
$ws3.Cells.Item(8, 4).Value = -939.2318524257653
$ws3.Cells.Item(9, 1).Value = "{'Age': np.int64(3), 'Student': np.int64(1)}"
$ws3.Cells.Item(9, 2).Value = 15934
$ws3.Cells.Item(9, 3).Value = 7355.426207414689
$ws3.Cells.Item(9, 4).Value = 642.0705062013676
$ws3.Cells.Item(10, 1).Value = "{'HDI': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(10, 2).Value = 23773
$ws3.Cells.Item(10, 3).Value = 7095.77558069452
$ws3.Cells.Item(10, 4).Value = 382.4198794811982
$ws3.Cells.Item(11, 1).Value = "{'GINI': np.int64(2), 'Student': np.int64(1)}"
$ws3.Cells.Item(11, 2).Value = 16046
$ws3.Cells.Item(11, 3).Value = 7226.922415426796
$ws3.Cells.Item(11, 4).Value = 513.5667142134744
$ws3.Cells.Item(12, 1).Value = "{'Gender': np.int64(1), 'FormalEducation': np.int64(1)}"
$ws3.Cells.Item(12, 2).Value = 16065
$ws3.Cells.Item(12, 3).Value = 4156.361432148196
$ws3.Cells.Item(12, 4).Value = -2556.994269065125
$ws3.Cells.Item(13, 1).Value = "{'Gender': np.int64(1), 'UndergradMajor': np.int64(2)}"
$ws3.Cells.Item(13, 2).Value = 21070
$ws3.Cells.Item(13, 3).Value = 2140.616974317039
$ws3.Cells.Item(13, 4).Value = -4572.738726896283
$ws3.Cells.Item(14, 1).Value = "{'UndergradMajor': np.int64(2), 'RaceEthnicity': np.int64(1)}"
$ws3.Cells.Item(14, 2).Value = 16908
$ws3.Cells.Item(14, 3).Value = 3943.333030186835
$ws3.Cells.Item(14, 4).Value = -2770.022671026486
$ws3.Cells.Item(15, 1).Value = "{'UndergradMajor': np.int64(2), 'Dependents': np.int64(2)}"
$ws3.Cells.Item(15, 2).Value = 15506
$ws3.Cells.Item(15, 3).Value = 6360.879966712541
$ws3.Cells.Item(15, 4).Value = -352.4757345007802
$ws3.Cells.Item(16, 1).Value = "{'UndergradMajor': np.int64(2), 'HDI': np.int64(1)}"
$ws3.Cells.Item(16, 2).Value = 17185
$ws3.Cells.Item(16, 3).Value = 3962.0157917935
$ws3.Cells.Item(16, 4).Value = -2751.339909419821
$ws3.Cells.Item(17, 1).Value = "{'Gender': np.int64(1), 'HoursComputer': np.int64(2)}"
$ws3.Cells.Item(17, 2).Value = 18591
$ws3.Cells.Item(17, 3).Value = 5759.491004058506
$ws3.Cells.Item(17, 4).Value = -953.8646971548151
$ws3.Cells.Item(18, 1).Value = "{'HoursComputer': np.int64(2), 'RaceEthnicity': np.int64(1)}"
$ws3.Cells.Item(18, 2).Value = 15402
$ws3.Cells.Item(18, 3).Value = 6060.471993983761
$ws3.Cells.Item(18, 4).Value = -652.8837072295601
$ws3.Cells.Item(19, 1).Value = "{'HoursComputer': np.int64(2), 'HDI': np.int64(1)}"
$ws3.Cells.Item(19, 2).Value = 15321
$ws3.Cells.Item(19, 3).Value = 6262.554982791322
$ws3.Cells.Item(19, 4).Value = -450.8007184219996
$ws3.Cells.Item(20, 1).Value = "{'Gender': np.int64(1), 'RaceEthnicity': np.int64(1)}"
$ws3.Cells.Item(20, 2).Value = 25910
$ws3.Cells.Item(20, 3).Value = 6471.012774252379
$ws3.Cells.Item(20, 4).Value = -242.3429269609423
$ws3.Cells.Item(21, 1).Value = "{'Hobby': np.int64(1)}"
$ws3.Cells.Item(21, 2).Value = 28842
$ws3.Cells.Item(21, 3).Value = 5516.21828962618
$ws3.Cells.Item(21, 4).Value = -1197.137411587141
$ws3.Cells.Item(22, 1).Value = "{'Student': np.int64(1)}"
$ws3.Cells.Item(22, 2).Value = 29526
$ws3.Cells.Item(22, 3).Value = 5617.758312088003
$ws3.Cells.Item(22, 4).Value = -1095.597389125319
$ws3.Cells.Item(23, 1).Value = "{'FormalEducation': np.int64(1)}"
$ws3.Cells.Item(23, 2).Value = 17131
$ws3.Cells.Item(23, 3).Value = 5336.323872777189
$ws3.Cells.Item(23, 4).Value = -1377.031828436133
$ws3.Cells.Item(24, 1).Value = "{'UndergradMajor': np.int64(2)}"
$ws3.Cells.Item(24, 2).Value = 22173
$ws3.Cells.Item(24, 3).Value = 3542.149042873403
$ws3.Cells.Item(24, 4).Value = -3171.206658339918
$ws3.Cells.Item(25, 1).Value = "{'HoursComputer': np.int64(2)}"
$ws3.Cells.Item(25, 2).Value = 19738
$ws3.Cells.Item(25, 3).Value = 5820.891848914629
$ws3.Cells.Item(25, 4).Value = -892.4638522986925
$ws3.Cells.Item(26, 1).Value = "{'Gender': np.int64(1)}"
$ws3.Cells.Item(26, 2).Value = 33253
$ws3.Cells.Item(26, 3).Value = 5920.984439082945
$ws3.Cells.Item(26, 4).Value = -792.3712621303766
$ws3.Cells.Item(27, 1).Value = "{'RaceEthnicity': np.int64(1)}"
$ws3.Cells.Item(27, 2).Value = 27379
$ws3.Cells.Item(27, 3).Value = 6509.936391650183
$ws3.Cells.Item(27, 4).Value = -203.4193095631381
$ws3.Cells.Item(28, 1).Value = "{'Dependents': np.int64(2)}"
$ws3.Cells.Item(28, 2).Value = 24167
$ws3.Cells.Item(28, 3).Value = 7072.162954162452
$ws3.Cells.Item(28, 4).Value = 358.8072529491301
$ws3.Cells.Item(29, 1).Value = "{'Age': np.int64(3)}"
$ws3.Cells.Item(29, 2).Value = 18401
$ws3.Cells.Item(29, 3).Value = 8048.308572796766
$ws3.Cells.Item(29, 4).Value = 1334.952871583445
$ws3.Cells.Item(30, 1).Value = "{'HDI': np.int64(1)}"
$ws3.Cells.Item(30, 2).Value = 27581
$ws3.Cells.Item(30, 3).Value = 8280.034330059932
$ws3.Cells.Item(30, 4).Value = 1566.678628846611
$ws3.Cells.Item(31, 1).Value = "{'GDP': np.int64(1)}"
$ws3.Cells.Item(31, 2).Value = 15752
$ws3.Cells.Item(31, 3).Value = 3295.830519890884
$ws3.Cells.Item(31, 4).Value = -3417.525181322438
$ws3.Cells.Item(32, 1).Value = "{'GINI': np.int64(2)}"
$ws3.Cells.Item(32, 2).Value = 19081
$ws3.Cells.Item(32, 3).Value = 9822.797811066226
$ws3.Cells.Item(32, 4).Value = 3109.442109852905
$ws3.Cells.Item(33, 1).Value = "{'Hobby': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(33, 2).Value = 23894
$ws3.Cells.Item(33, 3).Value = 4208.446115890501
$ws3.Cells.Item(33, 4).Value = -2504.90958532282
$ws3.Cells.Item(34, 1).Value = "{'UndergradMajor': np.int64(2), 'Hobby': np.int64(1)}"
$ws3.Cells.Item(34, 2).Value = 18239
$ws3.Cells.Item(34, 3).Value = 2972.364938069738
$ws3.Cells.Item(34, 4).Value = -3740.990763143584
$ws3.Cells.Item(35, 1).Value = "{'HoursComputer': np.int64(2), 'Hobby': np.int64(1)}"
$ws3.Cells.Item(35, 2).Value = 16363
$ws3.Cells.Item(35, 3).Value = 5602.899987220042
$ws3.Cells.Item(35, 4).Value = -1110.455713993279
$ws3.Cells.Item(36, 1).Value = "{'Gender': np.int64(1), 'Hobby': np.int64(1)}"
$ws3.Cells.Item(36, 2).Value = 27456
$ws3.Cells.Item(36, 3).Value = 5084.527001737077
$ws3.Cells.Item(36, 4).Value = -1628.828699476245
$ws3.Cells.Item(37, 1).Value = "{'RaceEthnicity': np.int64(1), 'Hobby': np.int64(1)}"
$ws3.Cells.Item(37, 2).Value = 22529
$ws3.Cells.Item(37, 3).Value = 5345.204113474854
$ws3.Cells.Item(37, 4).Value = -1368.151587738467
$ws3.Cells.Item(38, 1).Value = "{'Dependents': np.int64(2), 'Hobby': np.int64(1)}"
$ws3.Cells.Item(38, 2).Value = 19952
$ws3.Cells.Item(38, 3).Value = 5741.315462563783
$ws3.Cells.Item(38, 4).Value = -972.0402386495389
$ws3.Cells.Item(39, 1).Value = "{'Age': np.int64(3), 'Hobby': np.int64(1)}"
$ws3.Cells.Item(39, 2).Value = 15101
$ws3.Cells.Item(39, 3).Value = 7790.015822527435
$ws3.Cells.Item(39, 4).Value = 1076.660121314114
$ws3.Cells.Item(40, 1).Value = "{'Gender': np.int64(1), 'Dependents': np.int64(2)}"
$ws3.Cells.Item(40, 2).Value = 22568
$ws3.Cells.Item(40, 3).Value = 6660.009269424866
$ws3.Cells.Item(40, 4).Value = -53.34643178845545
$ws3.Cells.Item(41, 1).Value = "{'Gender': np.int64(1), 'Age': np.int64(3)}"
$ws3.Cells.Item(41, 2).Value = 17245
$ws3.Cells.Item(41, 3).Value = 8060.128048085273
$ws3.Cells.Item(41, 4).Value = 1346.772346871951
$ws3.Cells.Item(42, 1).Value = "{'Gender': np.int64(1), 'HDI': np.int64(1)}"
$ws3.Cells.Item(42, 2).Value = 25944
$ws3.Cells.Item(42, 3).Value = 7219.315451383841
$ws3.Cells.Item(42, 4).Value = 505.9597501705193
$ws3.Cells.Item(43, 1).Value = "{'GINI': np.int64(2), 'Gender': np.int64(1)}"
$ws3.Cells.Item(43, 2).Value = 17827
$ws3.Cells.Item(43, 3).Value = 9167.440434295811
$ws3.Cells.Item(43, 4).Value = 2454.08473308249
$ws3.Cells.Item(44, 1).Value = "{'RaceEthnicity': np.int64(1), 'Dependents': np.int64(2)}"
$ws3.Cells.Item(44, 2).Value = 18715
$ws3.Cells.Item(44, 3).Value = 5443.202083884486
$ws3.Cells.Item(44, 4).Value = -1270.153617328835
$ws3.Cells.Item(45, 1).Value = "{'RaceEthnicity': np.int64(1), 'HDI': np.int64(1)}"
$ws3.Cells.Item(45, 2).Value = 23146
$ws3.Cells.Item(45, 3).Value = 7669.000695259265
$ws3.Cells.Item(45, 4).Value = 955.6449940459433
$ws3.Cells.Item(46, 1).Value = "{'Dependents': np.int64(2), 'HDI': np.int64(1)}"
$ws3.Cells.Item(46, 2).Value = 18755
$ws3.Cells.Item(46, 3).Value = 8440.872413111409
$ws3.Cells.Item(46, 4).Value = 1727.516711898087
$ws3.Cells.Item(47, 1).Value = "{'HDI': np.int64(1), 'GDP': np.int64(1)}"
$ws3.Cells.Item(47, 2).Value = 15752
$ws3.Cells.Item(47, 3).Value = 3295.830519890884
$ws3.Cells.Item(47, 4).Value = -3417.525181322438
$ws3.Cells.Item(48, 1).Value = "{'Gender': np.int64(1), 'Hobby': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(48, 2).Value = 22758
$ws3.Cells.Item(48, 3).Value = 4081.320153739394
$ws3.Cells.Item(48, 4).Value = -2632.035547473928
$ws3.Cells.Item(49, 1).Value = "{'RaceEthnicity': np.int64(1), 'Hobby': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(49, 2).Value = 18956
$ws3.Cells.Item(49, 3).Value = 4707.214742762342
$ws3.Cells.Item(49, 4).Value = -2006.140958450979
$ws3.Cells.Item(50, 1).Value = "{'Dependents': np.int64(2), 'Hobby': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(50, 2).Value = 15886
$ws3.Cells.Item(50, 3).Value = 3923.495586716469
$ws3.Cells.Item(50, 4).Value = -2789.860114496852
$ws3.Cells.Item(51, 1).Value = "{'Hobby': np.int64(1), 'Student': np.int64(1), 'HDI': np.int64(1)}"
$ws3.Cells.Item(51, 2).Value = 19358
$ws3.Cells.Item(51, 3).Value = 5757.284204945694
$ws3.Cells.Item(51, 4).Value = -956.0714962676275
$ws3.Cells.Item(52, 1).Value = "{'Gender': np.int64(1), 'UndergradMajor': np.int64(2), 'Hobby': np.int64(1)}"
$ws3.Cells.Item(52, 2).Value = 17486
$ws3.Cells.Item(52, 3).Value = 1946.722347701861
$ws3.Cells.Item(52, 4).Value = -4766.633353511461
$ws3.Cells.Item(53, 1).Value = "{'Gender': np.int64(1), 'HoursComputer': np.int64(2), 'Hobby': np.int64(1)}"
$ws3.Cells.Item(53, 2).Value = 15555
$ws3.Cells.Item(53, 3).Value = 5831.184976203039
$ws3.Cells.Item(53, 4).Value = -882.1707250102827
$ws3.Cells.Item(54, 1).Value = "{'Gender': np.int64(1), 'RaceEthnicity': np.int64(1), 'Hobby': np.int64(1)}"
$ws3.Cells.Item(54, 2).Value = 21544
$ws3.Cells.Item(54, 3).Value = 5551.102665918384
$ws3.Cells.Item(54, 4).Value = -1162.253035294938
$ws3.Cells.Item(55, 1).Value = "{'Gender': np.int64(1), 'Dependents': np.int64(2), 'Hobby': np.int64(1)}"
$ws3.Cells.Item(55, 2).Value = 18868
$ws3.Cells.Item(55, 3).Value = 5634.634999361796
$ws3.Cells.Item(55, 4).Value = -1078.720701851526
$ws3.Cells.Item(56, 1).Value = "{'Gender': np.int64(1), 'Hobby': np.int64(1), 'HDI': np.int64(1)}"
$ws3.Cells.Item(56, 2).Value = 21574
$ws3.Cells.Item(56, 3).Value = 6412.706957617019
$ws3.Cells.Item(56, 4).Value = -300.6487435963027
$ws3.Cells.Item(57, 1).Value = "{'Dependents': np.int64(2), 'Hobby': np.int64(1), 'RaceEthnicity': np.int64(1)}"
$ws3.Cells.Item(57, 2).Value = 15532
$ws3.Cells.Item(57, 3).Value = 4235.626288481227
$ws3.Cells.Item(57, 4).Value = -2477.729412732095
$ws3.Cells.Item(58, 1).Value = "{'RaceEthnicity': np.int64(1), 'Hobby': np.int64(1), 'HDI': np.int64(1)}"
$ws3.Cells.Item(58, 2).Value = 19133
$ws3.Cells.Item(58, 3).Value = 6167.928095255636
$ws3.Cells.Item(58, 4).Value = -545.427605957686
$ws3.Cells.Item(59, 1).Value = "{'Dependents': np.int64(2), 'Hobby': np.int64(1), 'HDI': np.int64(1)}"
$ws3.Cells.Item(59, 2).Value = 15574
$ws3.Cells.Item(59, 3).Value = 7178.223410287983
$ws3.Cells.Item(59, 4).Value = 464.8677090746614
$ws3.Cells.Item(60, 1).Value = "{'Gender': np.int64(1), 'UndergradMajor': np.int64(2), 'Student': np.int64(1)}"
$ws3.Cells.Item(60, 2).Value = 17295
$ws3.Cells.Item(60, 3).Value = 1412.889200370405
$ws3.Cells.Item(60, 4).Value = -5300.466500842917
$ws3.Cells.Item(61, 1).Value = "{'Gender': np.int64(1), 'HoursComputer': np.int64(2), 'Student': np.int64(1)}"
$ws3.Cells.Item(61, 2).Value = 15737
$ws3.Cells.Item(61, 3).Value = 4570.176052280543
$ws3.Cells.Item(61, 4).Value = -2143.179648932779
$ws3.Cells.Item(62, 1).Value = "{'Gender': np.int64(1), 'RaceEthnicity': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(62, 2).Value = 22014
$ws3.Cells.Item(62, 3).Value = 5982.017555487582
$ws3.Cells.Item(62, 4).Value = -731.3381457257392
$ws3.Cells.Item(63, 1).Value = "{'Gender': np.int64(1), 'Dependents': np.int64(2), 'Student': np.int64(1)}"
$ws3.Cells.Item(63, 2).Value = 18145
$ws3.Cells.Item(63, 3).Value = 5331.204181537214
$ws3.Cells.Item(63, 4).Value = -1382.151519676107
$ws3.Cells.Item(64, 1).Value = "{'Gender': np.int64(1), 'HDI': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(64, 2).Value = 22370
$ws3.Cells.Item(64, 3).Value = 6176.689060932703
$ws3.Cells.Item(64, 4).Value = -536.6666402806186
$ws3.Cells.Item(65, 1).Value = "{'RaceEthnicity': np.int64(1), 'Dependents': np.int64(2), 'Student': np.int64(1)}"
$ws3.Cells.Item(65, 2).Value = 15301
$ws3.Cells.Item(65, 3).Value = 3964.682726108012
$ws3.Cells.Item(65, 4).Value = -2748.67297510531
$ws3.Cells.Item(66, 1).Value = "{'RaceEthnicity': np.int64(1), 'HDI': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(66, 2).Value = 20128
$ws3.Cells.Item(66, 3).Value = 6694.06933007496
$ws3.Cells.Item(66, 4).Value = -19.28637113836157
$ws3.Cells.Item(67, 1).Value = "{'Dependents': np.int64(2), 'HDI': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(67, 2).Value = 15601
$ws3.Cells.Item(67, 3).Value = 6881.083296062101
$ws3.Cells.Item(67, 4).Value = 167.7275948487795
$ws3.Cells.Item(68, 1).Value = "{'Gender': np.int64(1), 'UndergradMajor': np.int64(2), 'RaceEthnicity': np.int64(1)}"
$ws3.Cells.Item(68, 2).Value = 16152
$ws3.Cells.Item(68, 3).Value = 3274.466712898681
$ws3.Cells.Item(68, 4).Value = -3438.88898831464
$ws3.Cells.Item(69, 1).Value = "{'Gender': np.int64(1), 'UndergradMajor': np.int64(2), 'HDI': np.int64(1)}"
$ws3.Cells.Item(69, 2).Value = 16344
$ws3.Cells.Item(69, 3).Value = 2289.159933879372
$ws3.Cells.Item(69, 4).Value = -4424.195767333949
$ws3.Cells.Item(70, 1).Value = "{'Gender': np.int64(1), 'RaceEthnicity': np.int64(1), 'Dependents': np.int64(2)}"
$ws3.Cells.Item(70, 2).Value = 17563
$ws3.Cells.Item(70, 3).Value = 5093.2900472175
$ws3.Cells.Item(70, 4).Value = -1620.065653995822
$ws3.Cells.Item(71, 1).Value = "{'Gender': np.int64(1), 'RaceEthnicity': np.int64(1), 'HDI': np.int64(1)}"
$ws3.Cells.Item(71, 2).Value = 21900
$ws3.Cells.Item(71, 3).Value = 7566.667087338407
$ws3.Cells.Item(71, 4).Value = 853.3113861250858
$ws3.Cells.Item(72, 1).Value = "{'Gender': np.int64(1), 'Dependents': np.int64(2), 'HDI': np.int64(1)}"
$ws3.Cells.Item(72, 2).Value = 17456
$ws3.Cells.Item(72, 3).Value = 7663.615191867112
$ws3.Cells.Item(72, 4).Value = 950.2594906537906
$ws3.Cells.Item(73, 1).Value = "{'Dependents': np.int64(2), 'HDI': np.int64(1), 'RaceEthnicity': np.int64(1)}"
$ws3.Cells.Item(73, 2).Value = 15673
$ws3.Cells.Item(73, 3).Value = 6373.699265659056
$ws3.Cells.Item(73, 4).Value = -339.656435554266
$ws3.Cells.Item(74, 1).Value = "{'Gender': np.int64(1), 'RaceEthnicity': np.int64(1), 'Hobby': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(74, 2).Value = 18136
$ws3.Cells.Item(74, 3).Value = 5289.423209809513
$ws3.Cells.Item(74, 4).Value = -1423.932491403809
$ws3.Cells.Item(75, 1).Value = "{'Gender': np.int64(1), 'Dependents': np.int64(2), 'Hobby': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(75, 2).Value = 15007
$ws3.Cells.Item(75, 3).Value = 4003.268385218564
$ws3.Cells.Item(75, 4).Value = -2710.087315994757
$ws3.Cells.Item(76, 1).Value = "{'Gender': np.int64(1), 'Hobby': np.int64(1), 'Student': np.int64(1), 'HDI': np.int64(1)}"
$ws3.Cells.Item(76, 2).Value = 18438
$ws3.Cells.Item(76, 3).Value = 5619.548811679768
$ws3.Cells.Item(76, 4).Value = -1093.806889533554
$ws3.Cells.Item(77, 1).Value = "{'RaceEthnicity': np.int64(1), 'Hobby': np.int64(1), 'Student': np.int64(1), 'HDI': np.int64(1)}"
$ws3.Cells.Item(77, 2).Value = 16486
$ws3.Cells.Item(77, 3).Value = 5666.495872135472
$ws3.Cells.Item(77, 4).Value = -1046.859829077849
$ws3.Cells.Item(78, 1).Value = "{'Gender': np.int64(1), 'RaceEthnicity': np.int64(1), 'Hobby': np.int64(1), 'HDI': np.int64(1)}"
$ws3.Cells.Item(78, 2).Value = 18306
$ws3.Cells.Item(78, 3).Value = 6286.555506924179
$ws3.Cells.Item(78, 4).Value = -426.8001942891424
$ws3.Cells.Item(79, 1).Value = "{'Gender': np.int64(1), 'RaceEthnicity': np.int64(1), 'HDI': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(79, 2).Value = 19043
$ws3.Cells.Item(79, 3).Value = 7073.278421988637
$ws3.Cells.Item(79, 4).Value = 359.9227207753156
$ws3.Cells.Item(80, 1).Value = "{'Hobby': np.int64(1), 'Gender': np.int64(1), 'RaceEthnicity': np.int64(1), 'HDI': np.int64(1), 'Student': np.int64(1)}"
$ws3.Cells.Item(80, 2).Value = 15779
$ws3.Cells.Item(80, 3).Value = 6251.332734797039
$ws3.Cells.Item(80, 4).Value = -462.0229664162825
